$d = $word.ActiveDocument

# Hybrid bold + color (2C3E50) highlighting for quantitative impact metrics.
# Word's Font.Color takes a BGR-packed integer (R + G*256 + B*65536).
$highlightColor = 0x2C + (0x3E * 256) + (0x50 * 65536)

function Highlight-Metrics($ParagraphIndex, $Metrics) {
    $para = $d.Paragraphs($ParagraphIndex)
    $paraStart = $para.Range.Start
    $paraEnd = $para.Range.End

    $searchFrom = $paraStart
    foreach ($metric in $Metrics) {
        $scope = $d.Range($searchFrom, $paraEnd)
        $find = $scope.Find
        $find.Text = $metric
        $find.MatchCase = $true
        $found = $find.Execute()
        if ($found) {
            $scope.Font.Bold = 1
            $scope.Font.Color = $highlightColor
            $searchFrom = $scope.End
        }
    }
}

# Data Science & Political Analytics bullets (Siege Analytics)
Highlight-Metrics 10 @('23%', '64%')
Highlight-Metrics 12 @('±4.2%', '±2.1%', '71%', '87%')
Highlight-Metrics 13 @('73.5%', '$4.7M')
Highlight-Metrics 14 @('$2')

# Data Products Manager (Helm/Murmuration) bullet
Highlight-Metrics 24 @('57%')

# Key Achievements and Impact bullets
Highlight-Metrics 50 @('$4.9M')
Highlight-Metrics 51 @('23%')
Highlight-Metrics 53 @('12,847')
